$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-27 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-28 Thursday", 2)

$d.Content.Find.Execute("522×3=", $true, $false, $false, $false, $false, $true, 1, $false, "411×2=", 2)
$d.Content.Find.Execute("670×6=", $true, $false, $false, $false, $false, $true, 1, $false, "577×2=", 2)
$d.Content.Find.Execute("709×7=", $true, $false, $false, $false, $false, $true, 1, $false, "736×3=", 2)
$d.Content.Find.Execute("968×4=", $true, $false, $false, $false, $false, $true, 1, $false, "345×2=", 2)
$d.Content.Find.Execute("256×8=", $true, $false, $false, $false, $false, $true, 1, $false, "391×9=", 2)

$d.Content.Find.Execute("789×9=", $true, $false, $false, $false, $false, $true, 1, $false, "834×5=", 2)
$d.Content.Find.Execute("854×3=", $true, $false, $false, $false, $false, $true, 1, $false, "582×9=", 2)
$d.Content.Find.Execute("111×2=", $true, $false, $false, $false, $false, $true, 1, $false, "114×7=", 2)
$d.Content.Find.Execute("466×9=", $true, $false, $false, $false, $false, $true, 1, $false, "669×3=", 2)
$d.Content.Find.Execute("574×7=", $true, $false, $false, $false, $false, $true, 1, $false, "446×6=", 2)

$d.Content.Find.Execute("835×9=", $true, $false, $false, $false, $false, $true, 1, $false, "711×5=", 2)
$d.Content.Find.Execute("236×6=", $true, $false, $false, $false, $false, $true, 1, $false, "669×6=", 2)
$d.Content.Find.Execute("986×5=", $true, $false, $false, $false, $false, $true, 1, $false, "571×7=", 2)
$d.Content.Find.Execute("847×7=", $true, $false, $false, $false, $false, $true, 1, $false, "480×7=", 2)
$d.Content.Find.Execute("587×8=", $true, $false, $false, $false, $false, $true, 1, $false, "278×2=", 2)

$d.Content.Find.Execute("254×3=", $true, $false, $false, $false, $false, $true, 1, $false, "943×6=", 2)
$d.Content.Find.Execute("957×2=", $true, $false, $false, $false, $false, $true, 1, $false, "777×8=", 2)
$d.Content.Find.Execute("963×2=", $true, $false, $false, $false, $false, $true, 1, $false, "478×4=", 2)
$d.Content.Find.Execute("245×7=", $true, $false, $false, $false, $false, $true, 1, $false, "556×8=", 2)
$d.Content.Find.Execute("928×9=", $true, $false, $false, $false, $false, $true, 1, $false, "583×3=", 2)

$d.Content.Find.Execute("247×9=", $true, $false, $false, $false, $false, $true, 1, $false, "210×6=", 2)
$d.Content.Find.Execute("578×8=", $true, $false, $false, $false, $false, $true, 1, $false, "492×3=", 2)
$d.Content.Find.Execute("632×8=", $true, $false, $false, $false, $false, $true, 1, $false, "446×5=", 2)
$d.Content.Find.Execute("569×7=", $true, $false, $false, $false, $false, $true, 1, $false, "134×9=", 2)
$d.Content.Find.Execute("903×6=", $true, $false, $false, $false, $false, $true, 1, $false, "663×5=", 2)
